{"js": "const replacements = [\n  [\"513\u00f78=64, 1\", \"216\u00f76=36, 0\"],\n  [\"810\u00f74=202, 2\", \"618\u00f75=123, 3\"],\n  [\"876\u00f74=219, 0\", \"110\u00f78=13, 6\"],\n  [\"307\u00f79=34, 1\", \"403\u00f72=201, 1\"],\n  [\"238\u00f72=119, 0\", \"362\u00f76=60, 2\"],\n  [\"531\u00f79=59, 0\", \"546\u00f79=60, 6\"],\n  [\"734\u00f76=122, 2\", \"309\u00f78=38, 5\"],\n  [\"184\u00f73=61, 1\", \"314\u00f75=62, 4\"],\n  [\"981\u00f73=327, 0\", \"954\u00f77=136, 2\"],\n  [\"898\u00f76=149, 4\", \"465\u00f72=232, 1\"],\n  [\"170\u00f75=34, 0\", \"369\u00f79=41, 0\"],\n  [\"711\u00f73=237, 0\", \"524\u00f76=87, 2\"],\n  [\"550\u00f76=91, 4\", \"371\u00f73=123, 2\"],\n  [\"143\u00f79=15, 8\", \"441\u00f74=110, 1\"],\n  [\"924\u00f78=115, 4\", \"424\u00f79=47, 1\"],\n  [\"327\u00f79=36, 3\", \"451\u00f76=75, 1\"],\n  [\"247\u00f72=123, 1\", \"259\u00f72=129, 1\"],\n  [\"840\u00f79=93, 3\", \"586\u00f76=97, 4\"],\n  [\"513\u00f76=85, 3\", \"782\u00f77=111, 5\"],\n  [\"805\u00f79=89, 4\", \"373\u00f79=41, 4\"],\n  [\"200\u00f77=28, 4\", \"117\u00f75=23, 2\"],\n  [\"155\u00f77=22, 1\", \"179\u00f72=89, 1\"],\n  [\"400\u00f72=200, 0\", \"404\u00f72=202, 0\"],\n  [\"191\u00f72=95, 1\", \"931\u00f75=186, 1\"],\n  [\"627\u00f78=78, 3\", \"164\u00f78=20, 4\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"513\u00f78=64, 1\", \"216\u00f76=36, 0\"),\n    @(\"810\u00f74=202, 2\", \"618\u00f75=123, 3\"),\n    @(\"876\u00f74=219, 0\", \"110\u00f78=13, 6\"),\n    @(\"307\u00f79=34, 1\", \"403\u00f72=201, 1\"),\n    @(\"238\u00f72=119, 0\", \"362\u00f76=60, 2\"),\n    @(\"531\u00f79=59, 0\", \"546\u00f79=60, 6\"),\n    @(\"734\u00f76=122, 2\", \"309\u00f78=38, 5\"),\n    @(\"184\u00f73=61, 1\", \"314\u00f75=62, 4\"),\n    @(\"981\u00f73=327, 0\", \"954\u00f77=136, 2\"),\n    @(\"898\u00f76=149, 4\", \"465\u00f72=232, 1\"),\n    @(\"170\u00f75=34, 0\", \"369\u00f79=41, 0\"),\n    @(\"711\u00f73=237, 0\", \"524\u00f76=87, 2\"),\n    @(\"550\u00f76=91, 4\", \"371\u00f73=123, 2\"),\n    @(\"143\u00f79=15, 8\", \"441\u00f74=110, 1\"),\n    @(\"924\u00f78=115, 4\", \"424\u00f79=47, 1\"),\n    @(\"327\u00f79=36, 3\", \"451\u00f76=75, 1\"),\n    @(\"247\u00f72=123, 1\", \"259\u00f72=129, 1\"),\n    @(\"840\u00f79=93, 3\", \"586\u00f76=97, 4\"),\n    @(\"513\u00f76=85, 3\", \"782\u00f77=111, 5\"),\n    @(\"805\u00f79=89, 4\", \"373\u00f79=41, 4\"),\n    @(\"200\u00f77=28, 4\", \"117\u00f75=23, 2\"),\n    @(\"155\u00f77=22, 1\", \"179\u00f72=89, 1\"),\n    @(\"400\u00f72=200, 0\", \"404\u00f72=202, 0\"),\n    @(\"191\u00f72=95, 1\", \"931\u00f75=186, 1\"),\n    @(\"627\u00f78=78, 3\", \"164\u00f78=20, 4\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
